# Update "last updated" timestamp banner
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 13:22"

# Row 45 / 46: Catar overtakes Malasia in the "Casos totales" ranking, so the
# two countries swap places (their per-country stats are refreshed too).
$nameA45 = $ws.Range("A45").Text
$nameA46 = $ws.Range("A46").Text
$ws.Range("A45").Value = $nameA46
$ws.Range("A46").Value = $nameA45

# España (row 5)
$ws.Range("E5").Value = 98134
$ws.Range("G5").Value = 410
$ws.Range("H5").Value = 20453

# Iran (row 12)
$ws.Range("B12").Value = 82211
$ws.Range("C12").Value = 1343
$ws.Range("D12").Value = 57023
$ws.Range("E12").Value = 20070
$ws.Range("F12").Value = 3456
$ws.Range("G12").Value = 87
$ws.Range("H12").Value = 5118

# Portugal (row 19)
$ws.Range("B19").Value = 20206
$ws.Range("C19").Value = 521
$ws.Range("E19").Value = 18882
$ws.Range("F19").Value = 224
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 714

# Noruega (row 36)
$ws.Range("E36").Value = 6872
$ws.Range("F36").Value = 58
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 165

# Row 45 (now Catar) updated stats
$ws.Range("B45").Value = 5448
$ws.Range("C45").Value = 440
$ws.Range("D45").Value = 518
$ws.Range("E45").Value = 4922
$ws.Range("F45").Value = 37
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 8

# Row 46 (now Malasia) updated stats
$ws.Range("B46").Value = 5389
$ws.Range("C46").Value = 84
$ws.Range("D46").Value = 3197
$ws.Range("E46").Value = 2103
$ws.Range("F46").Value = 46
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 89

# Moldavia (row 60)
$ws.Range("B60").Value = 2351
$ws.Range("E60").Value = 1833

# Uzbekistan (row 70)
$ws.Range("D70").Value = 214
$ws.Range("E70").Value = 1276

# Bosnia y Herzegovina (row 76)
$ws.Range("B76").Value = 1285
$ws.Range("C76").Value = 17
$ws.Range("D76").Value = 347
$ws.Range("E76").Value = 890
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 48

# Republica de Macedonia (row 78)
$ws.Range("B78").Value = 1207
$ws.Range("C78").Value = 37
$ws.Range("D78").Value = 179
$ws.Range("E78").Value = 977
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 51

# San Marino (row 104)
$ws.Range("B104").Value = 461
$ws.Range("C104").Value = 6
$ws.Range("E104").Value = 362
$ws.Range("F104").Value = 4

# Malta (row 106)
$ws.Range("B106").Value = 427
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 118
$ws.Range("E106").Value = 306
$ws.Range("F106").Value = 2

# Brunei (row 132)
$ws.Range("B132").Value = 138
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 115
$ws.Range("E132").Value = 22

# Nepal (row 173)
$ws.Range("D173").Value = 4
$ws.Range("E173").Value = 27
